$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Remove the duplicated "Contact" row (old row 11); this shifts rows 12-22 up by one,
# shrinking the sheet from 22 to 21 rows.
$ws.Rows.Item(11).Delete()

# Version 5.0.0 -> 6.0.0
$ws.Range("B3").Value = "6.0.0"

# Date updated
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher now has a value
$ws.Range("B9").Value = "Alvearie Team"

# The remaining "Contact" row (old row 10) becomes the new "Jurisdiction" row
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"

# Case Sensitive now has a value of "true" (as literal text, not boolean).
# A direct .Value assignment of "true" gets auto-coerced to a boolean by Excel,
# so instead stage the text in a scratch cell with a quote-prefix, copy it, and
# paste-special (values only) into the target so the existing cell style/format
# of B14 is preserved and the value lands as plain text.
$ws.Range("Z1").Value = "'true"
$ws.Range("Z1").Copy()
$ws.Range("B14").PasteSpecial(-4163)
$ws.Range("Z1").Clear()
